$d = $word.ActiveDocument

# Locate the paragraph that holds the "{{ trial_court_name }}" placeholder.
$findRange = $d.Content
$found = $findRange.Find.Execute("trial_court_name", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'trial_court_name' placeholder text"
}
$targetPara = $findRange.Paragraphs(1)
$oldStart = $targetPara.Range.Start
$oldText = $targetPara.Range.Text

# Desired replacement text (splits trial_court_name into trial_court.address.county,
# surrounded by the gramStart/gramEnd proof-error markers Word adds around the edit).
$newText = "{{ trial_court.address.county }}"

# Build the new run/proofErr sequence and splice it in at the very start of the
# paragraph - inserting here keeps it as sibling content of the existing paragraph
# instead of spawning a new one.
$insertionPoint = $d.Range($oldStart, $oldStart)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>trial</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>_court</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.address.county</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xml)

# Remove the original "{{ trial_court_name }}" run content that now trails the
# freshly inserted replacement within the same paragraph.
$deleteStart = $oldStart + $newText.Length
$deleteEnd = $deleteStart + $oldText.Length
$d.Range($deleteStart, $deleteEnd).Delete()
